{"js": "// Update readme and play around with prompts\n// Rewrites the three Q&A answer paragraphs (surgeries, medications, allergies)\n// to reflect the new patient name (Jeanine Chase) and revised answer text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Line-break char used by Word for a soft return (<w:br/>) inside a run.\nconst BR = \"\\u000b\";\n\nconst replacements = [\n  {\n    match: \"Based on the medical data provided, the patient Deirdre Alimae Lum has had the following surgeries:\",\n    text:\n      \"Based on the provided medical data, Jeanine Chase has had the following surgeries:\" +\n      BR + BR +\n      \"* Mandibular and/or maxillary hardware removal (completed)\" +\n      BR + BR +\n      \"Note that the data only contains information about a few procedures/surgeries performed on Jeanine Chase. If you're looking for more information or a comprehensive list of her surgical history, I recommend reviewing the entire dataset or requesting additional medical records.\",\n  },\n  {\n    match: \"After reviewing the medical data, I found the following medications used by Deirdre Alimae Lum:\",\n    text:\n      \"Based on the medical data provided, Jeanine Chase has not been mentioned as a patient. The records appear to be related to Deirdre Alimae Lum and Robert Wayne Riley.\" +\n      BR + BR +\n      \"However, there is no mention of medications used by these individuals in the provided data. The records mainly consist of procedure notes, radiology reports, and laboratory test results.\",\n  },\n  {\n    match: \"There is no medical data that specifically indicates the patient has allergies.\",\n    text:\n      \"There is no medical data for Jeanine Chase in the provided text. The text appears to be a collection of medical records and procedures for individuals with different names (Deirdre Alimae Lum, Robert Wayne Riley, Harpreet Gill, Elizabeth Gomez, Kevin Cornwell, Donald M Sesso), but there is no mention or record of someone named Jeanine Chase.\",\n  },\n];\n\nfor (const r of replacements) {\n  const para = paragraphs.items.find((p) => p.text.startsWith(r.match));\n  if (para) {\n    para.insertText(r.text, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update readme and play around with prompts\n# Rewrites the three Q&A answer paragraphs (surgeries, medications, allergies)\n# to reflect the new patient name (Jeanine Chase) and revised answer text.\n\n$doc = $word.ActiveDocument\n\n$surgeriesOld = \"Based on the medical data provided, the patient*\"\n$surgeriesNew = \"Based on the provided medical data, Jeanine Chase has had the following surgeries:`v`v* Mandibular and/or maxillary hardware removal (completed)`v`vNote that the data only contains information about a few procedures/surgeries performed on Jeanine Chase. If you're looking for more information or a comprehensive list of her surgical history, I recommend reviewing the entire dataset or requesting additional medical records.\"\n\n$medsOld = \"After reviewing the medical data, I found the following medications*\"\n$medsNew = \"Based on the medical data provided, Jeanine Chase has not been mentioned as a patient. The records appear to be related to Deirdre Alimae Lum and Robert Wayne Riley.`v`vHowever, there is no mention of medications used by these individuals in the provided data. The records mainly consist of procedure notes, radiology reports, and laboratory test results.\"\n\n$allergiesOld = \"There is no medical data that specifically indicates the patient has allergies.*\"\n$allergiesNew = \"There is no medical data for Jeanine Chase in the provided text. The text appears to be a collection of medical records and procedures for individuals with different names (Deirdre Alimae Lum, Robert Wayne Riley, Harpreet Gill, Elizabeth Gomez, Kevin Cornwell, Donald M Sesso), but there is no mention or record of someone named Jeanine Chase.\"\n\nforeach ($p in $doc.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like $surgeriesOld) {\n        $p.Range.Text = $surgeriesNew\n    } elseif ($t -like $medsOld) {\n        $p.Range.Text = $medsNew\n    } elseif ($t -like $allergiesOld) {\n        $p.Range.Text = $allergiesNew\n    }\n}\n"}
